# Update scoresheets and scoring file
# TMFS_scoresheet: overwrite the old "select step 4" row (row 19) and the
# "mean" row (row 20) with a new set of "recode" steps (rows 19-24) that
# mean-center each tmfs item, push the final "mean" step down to row 25,
# and append a brand-new "tmfs_mean_center" step (row 26) driven by a
# TEXTJOIN() formula over the new *_ctr variable names.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# raw_vars (col A) for the six recode rows (19-24) - values already exist
# in the shared-string table (reused from rows 13-18), so fill these first.
$ws.Range("A19").Value = "tmfs_1_self_eval"
$ws.Range("A20").Value = "tmfs_2_ideal"
$ws.Range("A21").Value = "tmfs_3_interest"
$ws.Range("A22").Value = "tmfs_4_attitude"
$ws.Range("A23").Value = "tmfs_5_behavior"
$ws.Range("A24").Value = "tmfs_6_appearance"

# new_var (col B) - the new *_ctr variable names for each recode row.
$ws.Range("B19").Value = "tmfs_1_self_eval_ctr"
$ws.Range("B20").Value = "tmfs_2_ideal_ctr"
$ws.Range("B21").Value = "tmfs_3_interest_ctr"
$ws.Range("B22").Value = "tmfs_4_attitude_ctr"
$ws.Range("B23").Value = "tmfs_5_behavior_ctr"
$ws.Range("B24").Value = "tmfs_6_appearance_ctr"

# label (col C) is blank for all six recode rows.
$ws.Range("C19").Value = ""
$ws.Range("C20").Value = ""
$ws.Range("C21").Value = ""
$ws.Range("C22").Value = ""
$ws.Range("C23").Value = ""
$ws.Range("C24").Value = ""

# Row 26: brand-new "tmfs_mean_center" step.
$ws.Range("B26").Value = "tmfs_mean_center"
$ws.Range("C26").Value = "mean > 0 implies fem, mean < 0 implies masc"

# operation (col D) = "recode" for the six new rows.
$ws.Range("D19").Value = "recode"
$ws.Range("D20").Value = "recode"
$ws.Range("D21").Value = "recode"
$ws.Range("D22").Value = "recode"
$ws.Range("D23").Value = "recode"
$ws.Range("D24").Value = "recode"

# step (col E) = 3 for the six new recode rows.
$ws.Range("E19").Value = 3
$ws.Range("E20").Value = 3
$ws.Range("E21").Value = 3
$ws.Range("E22").Value = 3
$ws.Range("E23").Value = 3
$ws.Range("E24").Value = 3

# val_labs (col F) - the masculine/feminine scale description (-3..3).
$valLabs = "totally masculine =-3, masculine =-2, somewhat masculine =-1, netiher masculine nor feminine =0, somewhat feminine =1, feminine =2, totally feminine =3"
$ws.Range("F19").Value = $valLabs
$ws.Range("F20").Value = $valLabs
$ws.Range("F21").Value = $valLabs
$ws.Range("F22").Value = $valLabs
$ws.Range("F23").Value = $valLabs
$ws.Range("F24").Value = $valLabs

# new_vals (col G) - the 1..7 -> -3..3 recode mapping.
$newVals = "1=-3, 2 =-2, 3= -1, 4 = 0, 5 = 1, 6 = 2, 7= 3"
$ws.Range("G19").Value = $newVals
$ws.Range("G20").Value = $newVals
$ws.Range("G21").Value = $newVals
$ws.Range("G22").Value = $newVals
$ws.Range("G23").Value = $newVals
$ws.Range("G24").Value = $newVals

# if_condition / if_true_return / else_return / code (cols H-K) = "NA".
foreach ($r in 19..24) {
    $ws.Cells.Item($r, 8).Value = "NA"
    $ws.Cells.Item($r, 9).Value = "NA"
    $ws.Cells.Item($r, 10).Value = "NA"
    $ws.Cells.Item($r, 11).Value = "NA"
}

# Row 25: the original "tmfs_mean" step, now shifted down from row 20.
$ws.Range("A25").Value = "tmfs_1_self_eval, tmfs_2_ideal, tmfs_3_interest, tmfs_4_attitude, tmfs_5_behavior, tmfs_6_appearance"
$ws.Range("B25").Value = "tmfs_mean"
$ws.Range("C25").Value = "mean<4 implies masculinity, mean>4 implies femininity"
$ws.Range("D25").Value = "mean"
$ws.Range("E25").Value = 4
$ws.Range("F25").Value = "NA"
$ws.Range("G25").Value = "NA"
$ws.Range("H25").Value = "NA"
$ws.Range("I25").Value = "NA"
$ws.Range("J25").Value = "NA"
$ws.Range("K25").Value = "NA"

# Row 26 (continued): operation/step, plus a TEXTJOIN() formula in column A
# that builds the new raw_vars list from the *_ctr variable names.
$ws.Range("D26").Value = "mean"
$ws.Range("E26").Value = 4
$ws.Range("A26").Formula = '=TEXTJOIN(",",TRUE,B19:B24)'

# Move the selection to match the edited workbook's last-saved cursor position.
$ws.Range("G19").Select()
